# Fruta / hortaliza, semanal
# A new weekly price observation is inserted as row 96 (pushing the
# existing rows 96:120 down to 97:121), for
# "Feria Lagunitas de Puerto Montt" - Acelga.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 96; Excel copies the row-96
# (data-row) formatting from the row above it (row 95), which is what we
# want: default style everywhere except the date column (D) which keeps
# its date number-format style.
$ws.Rows("96:96").Insert()

# Populate the newly inserted row 96 with the new observation.
$ws.Range("A96").Value = 4
$ws.Range("B96").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C96").Value = "Los Lagos"
$ws.Range("D96").Value = 44511
$ws.Range("E96").Value = 10
$ws.Range("F96").Value = 100112009
$ws.Range("G96").Value = "Acelga"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 100
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 3000
$ws.Range("M96").Value = 3000
$ws.Range("N96").Value = "$/docena de atados (4 kilos)"
$ws.Range("O96").Value = "Región del Maule"
$ws.Range("P96").Value = 750
$ws.Range("Q96").Value = 4
$ws.Range("R96").Value = "Hortaliza"
